$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.259.23'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.113.50'
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.83'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.90'
$ws.Range("E6").Value = '  +0.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.51'
$ws.Range("E9").Value = '  +1.21%  '

$ws.Range("E10").Value = '  -0.89%  '

$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("E12").Value = '  -0.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.84'
$ws.Range("E13").Value = '  -1.15%  '

$ws.Range("E14").Value = '  -1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.632.28'
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.241.55'
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.11'
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.114.92'
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.65'
$ws.Range("E19").Value = '  +2.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '490.60'
$ws.Range("E20").Value = '  +0.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.89'
$ws.Range("E21").Value = '  +4.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").Value = '  -2.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.86'
$ws.Range("E23").Value = '  -0.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.18'
$ws.Range("E24").Value = '  -1.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  -2.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.59'
$ws.Range("E26").Value = '  +5.69%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.92'
$ws.Range("E28").Value = '  -1.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.37'
$ws.Range("E29").Value = '  -1.96%  '

$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.37'
$ws.Range("E31").Value = '  -2.04%  '

$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0938'
$ws.Range("E33").Value = '  -6.39%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("E35").Value = '  -0.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.971'
$ws.Range("E36").Value = '  -2.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '47.20'
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  -3.21%  '

$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("E40").Value = '  +1.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.48'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '392.80'
$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.804.16'
$ws.Range("E43").Value = '  -1.31%  '

$ws.Range("E44").Value = '  -7.67%  '

$ws.Range("E45").Value = '  -1.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.25'
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.01'
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("E49").Value = '  -0.60%  '

$ws.Range("E50").Value = '  -0.96%  '

$ws.Range("E51").Value = '  -2.07%  '
